# Fix the IN12 carrier BOM:
#  - Split the combined "1772-2080-ND 1772-1220-ND" Standoff row into two
#    separate rows (one per Digi-Key PN), adjusting quantities.
#  - Clean up the 470k resistor Digi-Key PN cell (was duplicated text).
#  - Append three new BOM lines: Pushbutton Cable, Pushbuttons, Feet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Split the Standoff row (currently row 12) into two rows ---
# Insert a new blank row at 13, pushing the old rows 13.. down by one.
$ws.Rows.Item(13).Insert()

# Row 12 keeps the same References/Value, but now only covers the
# "1772-2080-ND" parts, with an updated quantity.
$ws.Range("E12").Value = 8
$ws.Range("F12").Value = "1772-2080-ND "

# New row 13 covers the "1772-1220-ND" parts.
$ws.Range("A13").Value = $ws.Range("A12").Value()
$ws.Range("B13").Value = "Standoff"
$ws.Range("E13").Value = 16
$ws.Range("F13").Value = "1772-1220-ND"

# --- Clean up the duplicated Digi-Key PN text for the 470k resistors ---
# (That row was row 24 before the insert above; it is now row 25.)
$ws.Range("F25").Value = "541-470KCCT-ND"

# --- Append new BOM rows at the bottom ---
$ws.Range("A33").Value = "MK1001"
$ws.Range("B33").Value = "Pushbutton Cable"
$ws.Range("E33").Value = 1
$ws.Range("F33").Value = "SAM8931-ND"

$ws.Range("A34").Value = "MK1002 MK1003 MK1004 MK1005 MK1006"
$ws.Range("B34").Value = "Pushbuttons"
$ws.Range("E34").Value = 5
$ws.Range("F34").Value = "36-7600-ND"

$ws.Range("A35").Value = "MK1101"
$ws.Range("B35").Value = "Feet"
$ws.Range("E35").Value = 1
$ws.Range("F35").Value = "SJ5523-0-ND"

# Clear the lingering cell selection so the saved file doesn't pin the
# view to E2 (matches the author's final view state).
$ws.Range("A1").Select()
